$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Recorded By" values were generated by a set-based join whose ordering
# wasn't stable; re-normalize each multi-author cell in column G so the
# "System"/admin-ish entry that used to lead is moved to the end.
$lastRow = 157

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = @($val -split ", ")
        if ($parts.Length -gt 1) {
            $first = $parts[0]
            if ($first -eq "System" -or $first -eq "system" -or $first -eq "admin@admin.com") {
                $rest = $parts[1..($parts.Length - 1)]
                $rotated = $rest + @($first)
                $newVal = $rotated -join ", "
                $cell.Value = $newVal
            }
        }
    }
}
